$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before column N (shifts old N/O/P -> O/P/Q)
$ws.Columns("N").Insert()

# New column gets a custom width (stored width 10, i.e. ColumnWidth ~9.1667)
$ws.Columns("N").ColumnWidth = 9.1666666666667

# Make "Repayment Schedule" the active sheet/tab and update its selection
$ws.Activate()
$ws.Range("U6").Select()
